$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/savings-type"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# Remove the stray constraint text that had been duplicated onto the
# "Extension" row's Constraint(s) column (it correctly stays on the
# "Extension.extension" row, AI4).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
